$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value2
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 1.83 = 6631.29 pesos"), "✅ 1000 Bs = 1.79 = 6499.02 pesos"
$text = $text -replace [regex]::Escape("✅ 6631.29 pesos = 1.82 = 947.35 Bs"), "✅ 6499.02 pesos = 1.78 = 946.95 Bs"
$cell.Value = $text

# --- Sheet "tasas": update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 559.777
$ws2.Range("O10").Value = 3638
$ws2.Range("N12").Value = 3645
$ws2.Range("O12").Value = 531.1
